# Generate Report for Handoff
# Updates the "Priority" column (ht) and "Latest Handoff Datetime" values
# for the rows whose handoff has just completed, on both the zh-cn and
# de-de localization-status worksheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 14)

$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "ht"
    $ws.Range("H$r").Value = "2016-09-06 10:27:25"
}

$ws = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "ht"
    $ws.Range("H$r").Value = "2016-09-06 10:27:33"
}
